# "small edit to flow chart for groundfish in NE"
#
# 1) The cached "datetimeFigureOut" date field shown on the slide master and
#    every slide layout was re-stamped from 3/24/24 to 3/25/24 (PowerPoint
#    re-caches this text whenever the Date Placeholder is touched/saved).
# 2) The floating comment textbox ("TextBox 14") on the groundfish flow-chart
#    slide (slide 4) was repositioned.

$p = $ppt.ActivePresentation

# --- 1) Refresh the cached date-placeholder text -----------------------
$newDate = "3/25/24"

$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($L)
    $layoutShapes = $layout.Shapes
    for ($i = 1; $i -le $layoutShapes.Count; $i++) {
        $sh = $layoutShapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) Reposition the comment textbox on the groundfish flow chart ----
$slide = $p.Slides.Item(4)
$box = $slide.Shapes.Item("TextBox 14")
$box.Left = 965.30224609375
$box.Top = 301.9984436035156
